{"js": "// Apply the NATO/Norway edit described by the diff:\n//  Paragraph 1:\n//   - delete \"people \" (so \"5 million people and\" -> \"5 million and\")\n//   - replace \"a small country economically\" with \"its economy is small\"\n//  Paragraph 2:\n//   - delete \"deployment and \" (so \"Mobile Force) deployment and is exercised\" -> \"Mobile Force) is exercised\")\n//   - insert \" sometimes in very large scale\" right after \"exercised annually\"\n\nconst body = context.document.body;\n\n// --- Paragraph 1 edits ---\nconst hit1 = body.search(\"5 million people and\", { matchCase: true });\nhit1.load(\"items\");\nawait context.sync();\nif (hit1.items.length > 0) {\n  hit1.items[0].insertText(\"5 million and\", \"Replace\");\n  await context.sync();\n}\n\nconst hit2 = body.search(\"a small country economically\", { matchCase: true });\nhit2.load(\"items\");\nawait context.sync();\nif (hit2.items.length > 0) {\n  hit2.items[0].insertText(\"its economy is small\", \"Replace\");\n  await context.sync();\n}\n\n// --- Paragraph 2 edits ---\nconst hit3 = body.search(\"Mobile Force) deployment and is exercised\", { matchCase: true });\nhit3.load(\"items\");\nawait context.sync();\nif (hit3.items.length > 0) {\n  hit3.items[0].insertText(\"Mobile Force) is exercised\", \"Replace\");\n  await context.sync();\n}\n\nconst hit4 = body.search(\"exercised annually\", { matchCase: true });\nhit4.load(\"items\");\nawait context.sync();\nif (hit4.items.length > 0) {\n  hit4.items[0].insertText(\"exercised annually sometimes in very large scale\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Apply the NATO/Norway edit described by the diff:\n#  Paragraph 1:\n#   - delete \"people \" (so \"5 million people and\" -> \"5 million and\")\n#   - replace \"a small country economically\" with \"its economy is small\"\n#  Paragraph 2:\n#   - delete \"deployment and \" (so \"Mobile Force) deployment and is exercised\" -> \"Mobile Force) is exercised\")\n#   - insert \" sometimes in very large scale\" right after \"exercised annually\"\n\n$d = $word.ActiveDocument\n\n# --- Paragraph 1 edits ---\n$r1 = $d.Content\n$r1.Find.Execute(\"5 million people and\", $false, $false, $false, $false, $false, $true, 1, $false, \"5 million and\", 2) | Out-Null\n\n$r2 = $d.Content\n$r2.Find.Execute(\"a small country economically\", $false, $false, $false, $false, $false, $true, 1, $false, \"its economy is small\", 2) | Out-Null\n\n# --- Paragraph 2 edits ---\n$r3 = $d.Content\n$r3.Find.Execute(\"Mobile Force) deployment and is exercised\", $false, $false, $false, $false, $false, $true, 1, $false, \"Mobile Force) is exercised\", 2) | Out-Null\n\n$r4 = $d.Content\n$r4.Find.Execute(\"exercised annually\", $false, $false, $false, $false, $false, $true, 1, $false, \"exercised annually sometimes in very large scale\", 2) | Out-Null\n"}
